# Apply updated API call results to rows 6-19, columns B-E and J-M (duplicated "orig" columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    6  = @(2.890185992428072, 3.714905814451302, 1.515543661379421, 1.543688107481775)
    7  = @(3.947173685352226, 0.1045052961119888, 1.848995291837182, 1.902563758827545)
    8  = @(1.707089712754549, 2.999466478598363, 1.446321398849405, 1.780148400507346)
    9  = @(2.548265390879449, 5.161811764402577, 2.571210205683845, 1.924275155794865)
    10 = @(3.925161465394415, 4.573639312412173, 2.640929391352556, 2.066447829538879)
    11 = @(5.41570914763265,  8.690423321004072, 2.52100251750318,  2.160514947755103)
    12 = @(4.041327098464382, 7.104535720346792, 3.246296339093778, 2.405517533039561)
    13 = @(7.138470407773682, 6.223102706575908, 3.547903824342392, 2.500935234204537)
    14 = @(5.88510388946956,  8.686755909100066, 3.501288734535368, 2.657655480666553)
    15 = @(5.903278497842052, 9.326031465330262, 3.788507206940825, 2.875771655760958)
    16 = @(6.572391498800909, 9.855668632346275, 4.349832790633939, 2.717981976982983)
    17 = @(5.231883391341852, 1.673691006414161, 4.165709660639497, 2.591637474714854)
    18 = @(5.112999863523385, 6.10984356365371,  3.304065659169891, 2.632003991546259)
    19 = @(5.307787499498914, 4.979445338340163, 2.790295248693076, 2.786496694332275)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $b = $vals[0]
    $c = $vals[1]
    $d = $vals[2]
    $e = $vals[3]

    # gw_simul / gcpi_simul / cf1_simul / cf10_simul
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e

    # gw_simul_orig / gcpi_simul_orig / cf1_simul_orig / cf10_simul_orig
    $ws.Range("J$row").Value = $b
    $ws.Range("K$row").Value = $c
    $ws.Range("L$row").Value = $d
    $ws.Range("M$row").Value = $e
}
